$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the bold / centered / thin-border style on B1 first
$cell = $ws.Range("B1")
$cell.Borders.LineStyle = 1
$cell.Borders.Weight = 2
$cell.Font.Bold = $true
$cell.HorizontalAlignment = -4108
$cell.VerticalAlignment = -4160

# Re-use the exact same style on A2 via copy/paste-format so no
# duplicate intermediate style gets interned
$cell.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
